$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1861
$ws.Range("F3").Value = 23
$ws.Range("F5").Value = 58
$ws.Range("F8").Value = 657
$ws.Range("F9").Value = 81
$ws.Range("F10").Value = 476
$ws.Range("F11").Value = 796
$ws.Range("F12").Value = 1501
$ws.Range("F14").Value = 1459
$ws.Range("F16").Value = 1303
$ws.Range("F17").Value = 305
$ws.Range("F18").Value = 1610
$ws.Range("F19").Value = 782
$ws.Range("F20").Value = 1033
$ws.Range("F21").Value = 336
$ws.Range("F24").Value = 1461
$ws.Range("F26").Value = 143
$ws.Range("F27").Value = 802
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 990
$ws.Range("F32").Value = 36
$ws.Range("F33").Value = 558
$ws.Range("F34").Value = 1329
$ws.Range("F35").Value = 1070
$ws.Range("F37").Value = 1082
$ws.Range("F38").Value = 37
$ws.Range("F39").Value = 114
$ws.Range("F40").Value = 49
$ws.Range("F42").Value = 1618
$ws.Range("F44").Value = 56
$ws.Range("F45").Value = 797

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 42
$ws.Range("F7").Value = 1459
$ws.Range("F12").Value = 397
$ws.Range("F16").Value = 68
$ws.Range("F19").Value = 440
$ws.Range("F20").Value = 21
$ws.Range("F22").Value = 294
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 18
$ws.Range("F27").Value = 185
$ws.Range("F30").Value = 159
$ws.Range("F35").Value = 172
$ws.Range("F39").Value = 48
$ws.Range("F40").Value = 48
$ws.Range("F41").Value = 3
$ws.Range("F42").Value = 132
$ws.Range("F43").Value = 57

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 251
$ws.Range("F5").Value = 2821
$ws.Range("F6").Value = 4551
$ws.Range("F7").Value = 125
$ws.Range("F8").Value = 7
$ws.Range("F9").Value = 546
$ws.Range("F10").Value = 673
$ws.Range("F12").Value = 258
$ws.Range("F13").Value = 877
$ws.Range("F14").Value = 226
$ws.Range("F15").Value = 513

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1861
$ws.Range("F3").Value = 251
$ws.Range("F4").Value = 2821
$ws.Range("F5").Value = 4551
$ws.Range("F6").Value = 673
$ws.Range("F7").Value = 58
$ws.Range("F8").Value = 258
$ws.Range("F9").Value = 258
$ws.Range("F10").Value = 877
$ws.Range("F11").Value = 877
$ws.Range("F14").Value = 1459
$ws.Range("F15").Value = 476
$ws.Range("F16").Value = 796
$ws.Range("F19").Value = 1501
$ws.Range("F21").Value = 1459
$ws.Range("F22").Value = 1303
$ws.Range("F24").Value = 68
$ws.Range("F25").Value = 1610
$ws.Range("F26").Value = 782
$ws.Range("F27").Value = 1033
$ws.Range("F28").Value = 336
$ws.Range("F29").Value = 513
$ws.Range("F30").Value = 513
$ws.Range("F31").Value = 440
$ws.Range("F32").Value = 1461
$ws.Range("F33").Value = 143
$ws.Range("F34").Value = 802
$ws.Range("F37").Value = 294
$ws.Range("F38").Value = 990
$ws.Range("F39").Value = 36
$ws.Range("F40").Value = 1070
$ws.Range("F42").Value = 1082
$ws.Range("F44").Value = 114
$ws.Range("F47").Value = 1618
$ws.Range("F49").Value = 797
$ws.Range("F50").Value = 48
$ws.Range("F52").Value = 3
$ws.Range("F53").Value = 57
